$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 130.33333
$ws.Range("I9").Value = 150
$ws.Range("J9").Value = 120.5
$ws.Range("K9").Value = 150
$ws.Range("L9").Value = 120.5
$ws.Range("M9").Value = 19
$ws.Range("N9").Value = -458.5

$ws.Range("H40").Value = 2045.1818
$ws.Range("I40").Value = 1600
$ws.Range("J40").Value = 2212.125
$ws.Range("K40").Value = 1600
$ws.Range("L40").Value = 2212.125
$ws.Range("M40").Value = -1425
$ws.Range("N40").Value = -2562.125

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H97").Value = 1000
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992

$ws.Range("H98").Value = 771.125
$ws.Range("I98").Value = 451.85715
$ws.Range("K98").Value = 451.85715
$ws.Range("M98").Value = 1046.14285

$ws.Range("H116").Value = 4747.75
$ws.Range("I116").Value = 2997
$ws.Range("J116").Value = 6498.5
$ws.Range("K116").Value = 2997
$ws.Range("L116").Value = 6498.5
$ws.Range("M116").Value = 445
$ws.Range("N116").Value = -13382.5

$ws.Range("H122").Value = 771.125
$ws.Range("I122").Value = 451.85715
$ws.Range("K122").Value = 1355.57145
$ws.Range("M122").Value = 1094.42855

$ws.Range("H132").Value = 2169.5833
$ws.Range("I132").Value = 2169.5833
$ws.Range("K132").Value = 6508.749899999999
$ws.Range("M132").Value = -3978.749899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1939.4
$ws.Range("I61").Value = 1821.5555
$ws.Range("K61").Value = 1821.5555
$ws.Range("M61").Value = -1609.5555

$ws.Range("H74").Value = 1406.0834
$ws.Range("I74").Value = 1406.0834
$ws.Range("K74").Value = 1406.0834
$ws.Range("M74").Value = -532.0834

$ws.Range("H77").Value = 1406.0834
$ws.Range("I77").Value = 1406.0834
$ws.Range("K77").Value = 7030.416999999999
$ws.Range("M77").Value = -2662.416999999999

$ws.Range("H88").Value = 3524.5454
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3524.5454
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3524.5454
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4336.5454

$ws.Range("H91").Value = 3524.5454
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3524.5454
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3524.5454
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6332.5454

$ws.Range("H136").Value = 1939.4
$ws.Range("I136").Value = 1821.5555
$ws.Range("K136").Value = 5464.666499999999
$ws.Range("M136").Value = -2914.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 30251.096
$ws.Range("I82").Value = 15405.375
$ws.Range("J82").Value = 39386.92
$ws.Range("K82").Value = 15405.375
$ws.Range("L82").Value = 39386.92
$ws.Range("M82").Value = -15022.375
$ws.Range("N82").Value = -40152.92

$ws.Range("H85").Value = 30251.096
$ws.Range("I85").Value = 15405.375
$ws.Range("J85").Value = 39386.92
$ws.Range("K85").Value = 15405.375
$ws.Range("L85").Value = 39386.92
$ws.Range("M85").Value = -14079.375
$ws.Range("N85").Value = -42038.92

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3081.0908
$ws.Range("J31").Value = 5308.3335
$ws.Range("L31").Value = 5308.3335
$ws.Range("N31").Value = -5898.3335

$ws.Range("H34").Value = 3081.0908
$ws.Range("J34").Value = 5308.3335
$ws.Range("L34").Value = 5308.3335
$ws.Range("N34").Value = -5712.3335

$ws.Range("H58").Value = 4337.8237
$ws.Range("I58").Value = 3529.0833
$ws.Range("K58").Value = 3529.0833
$ws.Range("M58").Value = -3326.0833

$ws.Range("H88").Value = 7241.1665
$ws.Range("J88").Value = 8321.556
$ws.Range("L88").Value = 8321.556
$ws.Range("N88").Value = -9133.556

$ws.Range("H91").Value = 7241.1665
$ws.Range("J91").Value = 8321.556
$ws.Range("L91").Value = 8321.556
$ws.Range("N91").Value = -11129.556

$ws.Range("H136").Value = 4337.8237
$ws.Range("I136").Value = 3529.0833
$ws.Range("K136").Value = 10587.2499
$ws.Range("M136").Value = -8037.249899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 946.8570999999999
$ws.Range("I107").Value = 1348.5
$ws.Range("J107").Value = 786.2
$ws.Range("K107").Value = 4045.5
$ws.Range("L107").Value = 2358.6
$ws.Range("M107").Value = -2125.5
$ws.Range("N107").Value = -6198.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3006.3333
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3006.3333
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3006.3333
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5002.3333

$ws.Range("H83").Value = 3006.3333
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3006.3333
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 15031.6665
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -25015.6665

$ws.Range("H126").Value = 16412.375
$ws.Range("I126").Value = 14471.286
$ws.Range("K126").Value = 43413.858
$ws.Range("M126").Value = -40943.858

$ws.Range("H132").Value = 4913.25
$ws.Range("I132").Value = 4879.6665
$ws.Range("K132").Value = 14638.9995
$ws.Range("M132").Value = -12108.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 54122.25
$ws.Range("I7").Value = 54122.25
$ws.Range("K7").Value = 54122.25
$ws.Range("M7").Value = -54010.25

$ws.Range("H126").Value = 54122.25
$ws.Range("I126").Value = 54122.25
$ws.Range("K126").Value = 162366.75
$ws.Range("M126").Value = -159896.75

$ws.Range("H132").Value = 3854.25
$ws.Range("I132").Value = 1180.1428
$ws.Range("J132").Value = 7598
$ws.Range("K132").Value = 3540.4284
$ws.Range("L132").Value = 22794
$ws.Range("M132").Value = -1010.4284
$ws.Range("N132").Value = -27854

$ws.Range("H136").Value = 1787.7646
$ws.Range("I136").Value = 1664.5834
$ws.Range("J136").Value = 2083.4
$ws.Range("K136").Value = 4993.7502
$ws.Range("L136").Value = 6250.200000000001
$ws.Range("M136").Value = -2443.7502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1198.5
$ws.Range("I132").Value = 1198.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3595.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1065.5
$ws.Range("N132").ClearContents()
